$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (single decimal point,
# no thousands-separators) need NumberFormat forced to Text first, otherwise
# Excel auto-converts the entry to a numeric value and formatting such as
# trailing zeros ("1.000" -> 1) or leading zeros in small decimals would be lost.
$textCells = @(
    'D4',
    'D5',
    'D6',
    'D9',
    'D10',
    'D11',
    'D13',
    'D14',
    'D15',
    'D16',
    'D19',
    'D20',
    'D21',
    'D22',
    'D25',
    'D26',
    'D27',
    'D28',
    'D29',
    'D30',
    'D31',
    'D32',
    'D34',
    'D35',
    'D36',
    'D38',
    'D39',
    'D41',
    'D42',
    'D43',
    'D44',
    'D46',
    'D47',
    'D48',
    'D49',
    'D51'
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '24.968.54'
$ws.Range('E2').Value = '  +2.07%  '
$ws.Range('D3').Value = '1.700.63'
$ws.Range('E3').Value = '  +0.83%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '315.72'
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.23%  '
$ws.Range('E7').Value = '  +1.71%  '
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('D9').Value = '1.468'
$ws.Range('E9').Value = '  -1.38%  '
$ws.Range('D10').Value = '52.94'
$ws.Range('E10').Value = '  +1.05%  '
$ws.Range('D11').Value = '1.000'
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('E12').Value = '  +0.27%  '
$ws.Range('D13').Value = '25.95'
$ws.Range('E13').Value = '  -2.37%  '
$ws.Range('D14').Value = '7.458'
$ws.Range('E14').Value = '  +0.04%  '
$ws.Range('D15').Value = '0.00001351'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').Value = '7.956'
$ws.Range('E16').Value = '  -2.40%  '
$ws.Range('D17').Value = '1.705.90'
$ws.Range('E17').Value = '  +1.67%  '
$ws.Range('E18').Value = '  -2.20%  '
$ws.Range('D19').Value = '0.07189'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').Value = '20.66'
$ws.Range('E20').Value = '  +2.00%  '
$ws.Range('D21').Value = '7.322'
$ws.Range('E21').Value = '  +0.48%  '
$ws.Range('D22').Value = '1.001'
$ws.Range('E23').Value = '  +0.52%  '
$ws.Range('D24').Value = '24.973.88'
$ws.Range('E24').Value = '  +2.12%  '
$ws.Range('D25').Value = '2.365'
$ws.Range('E25').Value = '  +1.17%  '
$ws.Range('D26').Value = '2.948'
$ws.Range('E26').Value = '  -2.68%  '
$ws.Range('D27').Value = '23.72'
$ws.Range('E27').Value = '  +4.83%  '
$ws.Range('D28').Value = '6.242'
$ws.Range('E28').Value = '  +16.82%  '
$ws.Range('D29').Value = '162.48'
$ws.Range('E29').Value = '  -3.05%  '
$ws.Range('D30').Value = '150.70'
$ws.Range('E30').Value = '  +9.00%  '
$ws.Range('D31').Value = '8.343'
$ws.Range('E31').Value = '  -1.10%  '
$ws.Range('D32').Value = '2.631'
$ws.Range('E32').Value = '  +26.20%  '
$ws.Range('D33').Value = '1.896.09'
$ws.Range('E33').Value = '  +1.80%  '
$ws.Range('D34').Value = '0.08560'
$ws.Range('E34').Value = '  -2.08%  '
$ws.Range('D35').Value = '0.03148'
$ws.Range('E35').Value = '  +4.44%  '
$ws.Range('D36').Value = '7.157'
$ws.Range('E36').Value = '  -1.54%  '
$ws.Range('E37').Value = '  -1.07%  '
$ws.Range('D38').Value = '0.2866'
$ws.Range('E38').Value = '  +2.72%  '
$ws.Range('D39').Value = '0.09583'
$ws.Range('E39').Value = '  +4.83%  '
$ws.Range('E40').Value = '  +0.27%  '
$ws.Range('D41').Value = '0.8241'
$ws.Range('E41').Value = '  +2.70%  '
$ws.Range('D42').Value = '14.00'
$ws.Range('E42').Value = '  -1.09%  '
$ws.Range('D43').Value = '1.484'
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('D44').Value = '17.16'
$ws.Range('E44').Value = '  -2.77%  '
$ws.Range('E45').Value = '  +1.07%  '
$ws.Range('D46').Value = '0.7383'
$ws.Range('E46').Value = '  +1.80%  '
$ws.Range('D47').Value = '4.242'
$ws.Range('E47').Value = '  -0.50%  '
$ws.Range('D48').Value = '1.391'
$ws.Range('E48').Value = '  -1.74%  '
$ws.Range('D49').Value = '0.08797'
$ws.Range('E49').Value = '  +8.86%  '
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').Value = '139.20'
$ws.Range('E51').Value = '  -0.12%  '
